$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 3 (Jaymin Patel, 160130107024) and shift row 4 up.
$ws.Rows.Item(3).Delete()
